$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume deltas), GitHub Actions run.

$ws.Cells.Item(2, 4).Value = '66.623.31'
$ws.Cells.Item(2, 5).Value = '  -1.45%  '
$ws.Cells.Item(3, 4).Value = '3.796.01'
$ws.Cells.Item(3, 5).Value = '  -1.74%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.19%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '434.00'
$ws.Cells.Item(5, 5).Value = '  +4.62%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '139.48'
$ws.Cells.Item(6, 5).Value = '  +4.34%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.623'
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 5).Value = '  -1.87%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.153'
$ws.Cells.Item(10, 5).Value = '  -14.35%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0000318'
$ws.Cells.Item(11, 5).Value = '  -18.79%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '42.59'
$ws.Cells.Item(12, 5).Value = '  +2.49%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '10.46'
$ws.Cells.Item(13, 5).Value = '  +3.40%  '
$ws.Cells.Item(14, 4).Value = '4.407.28'
$ws.Cells.Item(14, 5).Value = '  -0.98%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '15.06'
$ws.Cells.Item(15, 5).Value = '  +1.44%  '
$ws.Cells.Item(17, 4).Value = '3.797.95'
$ws.Cells.Item(17, 5).Value = '  -1.58%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '19.96'
$ws.Cells.Item(18, 5).Value = '  +0.98%  '
$ws.Cells.Item(19, 5).Value = '  +3.57%  '
$ws.Cells.Item(20, 4).Value = '66.793.33'
$ws.Cells.Item(20, 5).Value = '  -1.25%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '411.56'
$ws.Cells.Item(21, 5).Value = '  -2.11%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '14.70'
$ws.Cells.Item(22, 5).Value = '  -2.50%  '
$ws.Cells.Item(23, 5).Value = '  +4.94%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '85.35'
$ws.Cells.Item(24, 5).Value = '  -2.19%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '36.90'
$ws.Cells.Item(25, 5).Value = '  -0.26%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.33'
$ws.Cells.Item(26, 5).Value = '  +4.53%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.78'
$ws.Cells.Item(27, 5).Value = '  +33.76%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.62'
$ws.Cells.Item(28, 5).Value = '  -2.02%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.82'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.136'
$ws.Cells.Item(30, 5).Value = '  +10.70%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '13.87'
$ws.Cells.Item(31, 5).Value = '  +10.15%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '714.78'
$ws.Cells.Item(32, 5).Value = '  +2.45%  '
$ws.Cells.Item(33, 5).Value = '  +2.30%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '41.60'
$ws.Cells.Item(34, 5).Value = '  +5.74%  '
$ws.Cells.Item(35, 5).Value = '  -0.08%  '
$ws.Cells.Item(36, 2).Value = 'NEARProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.67'
$ws.Cells.Item(36, 5).Value = '  +27.19%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.150'
$ws.Cells.Item(37, 5).Value = '  -3.48%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '55.70'
$ws.Cells.Item(38, 5).Value = '  -0.02%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0472'
$ws.Cells.Item(39, 5).Value = '  +0.99%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.78'
$ws.Cells.Item(40, 5).Value = '  +40.94%  '
$ws.Cells.Item(41, 2).Value = 'ThetaToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.94'
$ws.Cells.Item(41, 5).Value = '  -5.10%  '
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(42, 4).Value = '0.0₃0697'
$ws.Cells.Item(42, 5).Value = '  -13.09%  '
$ws.Cells.Item(43, 5).Value = '  +2.52%  '
$ws.Cells.Item(44, 5).Value = '  +0.73%  '
$ws.Cells.Item(45, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.22'
$ws.Cells.Item(45, 5).Value = '  +1.14%  '
$ws.Cells.Item(46, 2).Value = 'TheGraph'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.324'
$ws.Cells.Item(46, 5).Value = '  +8.32%  '
$ws.Cells.Item(47, 5).Value = '  -0.97%  '
$ws.Cells.Item(48, 5).Value = '  +3.20%  '
$ws.Cells.Item(49, 5).Value = '  -2.75%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '142.14'
$ws.Cells.Item(50, 5).Value = '  -4.26%  '
$ws.Cells.Item(51, 5).Value = '  -2.13%  '
